# Modification de la gestion des quantités de repas offerts pour les
# artistes et les bénévoles : on fusionne les deux colonnes de dates
# ("20/06/2025" et "21/06/2025") en une seule colonne "Quantité",
# en conservant les quantités qui étaient renseignées dans la seconde
# colonne de date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C currently holds "20/06/2025" (mostly empty) and column D holds
# "21/06/2025" (the actual quantities). Move D's values+formatting into C
# for the data rows (2-4), then drop column D. Row 1 (the header) is
# handled separately below since it becomes the new "Quantité" label.
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 4).Copy($ws.Cells.Item($r, 3))
}

# Rename the header of column C (still the "20/06/2025" shared string) to
# "Quantité" - this also renames the underlying table column.
$ws.Cells.Item(1, 3).Value2 = "Quantité"

# Drop the now-redundant last table column ("21/06/2025") and its
# worksheet column.
$lo = $ws.ListObjects.Item(1)
$cols = $lo.ListColumns
$cols.Item(4).Delete()
$ws.Columns.Item(4).Delete()

# Row 4 used to carry an explicit custom row format (tied to the old,
# now-unused "0 h" number format on the old column C). Clear it so the
# row goes back to plain per-cell formatting...
$ws.Rows.Item(4).ClearFormats()

# ...then restore A4/B4/C4's alignment + font so they still look like the
# rest of the "Artiste"/"Bénévole" rows (vertical-center, and
# horizontal-center for the category/quantity cells).
$ws.Cells.Item(4, 1).VerticalAlignment = -4108
$ws.Cells.Item(4, 1).Font.Name = "Calibri"
$ws.Cells.Item(4, 2).VerticalAlignment = -4108
$ws.Cells.Item(4, 2).HorizontalAlignment = -4108
$ws.Cells.Item(4, 2).Font.Name = "Calibri"
$ws.Cells.Item(4, 3).VerticalAlignment = -4108
$ws.Cells.Item(4, 3).HorizontalAlignment = -4108
$ws.Cells.Item(4, 3).Font.Name = "Calibri"

# Match the saved selection state.
$ws.Range("C2").Select()
